# The "2024" sheet keeps a running log of monthly transaction notes.
# A new September entry ("bal axisbank" at 2024-09-09 12:19:33) was logged
# at the top of the September_Details/September_Date list (columns R/S).
# Inserting a whole row at row 35 pushes that list - and everything below it
# in columns A/P/Q/R/S (the August list and the trailing "Broadband" label) -
# down by one row, growing the sheet from A1:Y117 to A1:Y118, which matches
# the recorded diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Shift row 35 (and everything below it) down by one row.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row with the new September log entry.
$ws.Range("R35").Value = "bal axisbank"
$ws.Range("S35").Value = "2024-09-09 12:19:33"
